# SG_Ag_compendium/Data/beef_sector.xlsx - "more updates for dairy/beef"
# Update the 2023 beef sector figures on Sheet1 and leave the selection
# where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of beef cows (2023): 394709 -> 394698
$ws.Range("D2").Value = 394698

# Number of holdings with beef cows (2023): 8134 -> 8127
$ws.Range("D3").Value = 8127

# Move the selection to C5, matching the saved cursor position.
$ws.Range("C5").Select()
